$wb = $excel.ActiveWorkbook

# --- Sheet 1: _set_CASES ---
# Header rename: c_Name -> cases_Name  (data below unchanged: "reference")
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "cases_Name"
$ws1.Range("A2").Select()

# --- Sheet 2: _set_PRODUCTS ---
# Header rename: p_Name -> products_Name, plus new column B header "products_Aggregation"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "products_Name"
$ws2.Range("A1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)  # xlPasteFormats - match header style of A1
$ws2.Range("B1").Value = "products_Aggregation"
$ws2.Range("A2:A4").Select()

# --- Sheet 3: _set_TECHNOLOGIES ---
# Header rename: t_Name -> technologies_Name, plus new column B header "technologies_Aggregation"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "technologies_Name"
$ws3.Range("A1").Copy()
$ws3.Range("B1").PasteSpecial(-4122)  # xlPasteFormats - match header style of A1
$ws3.Range("B1").Value = "technologies_Aggregation"
$ws3.Activate()
$ws3.Range("C10").Select()
